$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Neo4j/Web data file names (columns D and E) ---
$ws.Range("D2").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_WebData.xlsx"
$ws.Range("D3").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_WebData.xlsx"
$ws.Range("D4").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_WebData.xlsx"

# --- Update the Cypher queries (columns B and C) to filter on 'Illumina MiSeq' instead of 'DNBSEQ-G400' ---

# Row 2: Participants query (column B) + Files summary/count query (column C)
$ws.Range("B2").Value = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina MiSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p, s, collect(distinct samp.sample_id) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY ``Participant ID```LIMIT 100"

$filesSummaryQuery = "MATCH (f:file)`nMatch (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina MiSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,f, s, collect(distinct samp.sample_id) as samp`nRETURN`ncount(distinct s) AS Studies,`ncount(distinct p) AS Participants,`ncount(distinct samp) AS Samples,`ncount(distinct f) AS Files"

$ws.Range("C2").Value = $filesSummaryQuery

# Row 3: Samples query (column B)
$ws.Range("B3").Value = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina MiSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN  `n coalesce(samp.sample_id, '') as ``Sample ID``,`n coalesce(p.participant_id,'') as ``Participant ID``,`n coalesce(s.study_name, '') as ``Study Name``,`n coalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(samp.sample_tumor_status,'') as ``Tumor``,`ncoalesce(samp.sample_type,'') as ``Analyte Type```nORDER By samp.sample_id LIMIT 100"

$ws.Range("C3").Value = $filesSummaryQuery

# Row 4: Files query (column B)
$ws.Range("B4").Value = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina MiSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN `n    coalesce(f.file_name, '') as ``File Name``,`n    coalesce(s.study_name, '') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(p.participant_id,'') as ``Participant ID``,`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(f.file_type, '') as ``File Type```nORDER By f.file_name LIMIT 100"

$ws.Range("C4").Value = $filesSummaryQuery

# --- Column widths for D and E ---
# (Target stored widths are 88.85546875 / 87.140625; the host's column-width
#  model quantizes to multiples of 1/6, so these inputs land on the closest
#  representable stored widths: 88.833333.. / 87.166666..)
$ws.Columns.Item(4).ColumnWidth = 88.0
$ws.Columns.Item(5).ColumnWidth = 86.333333

# --- Move the active selection to D4 ---
$ws.Range("D4").Select()
